$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 13538.3
$ws.Range("J51").Value = 8376
$ws.Range("L51").Value = 8376
$ws.Range("N51").Value = -9344

$ws.Range("H52").Value = 298.86957
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 298.86957
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 896.60871
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -1216.60871

$ws.Range("H88").Value = 15194608
$ws.Range("I88").Value = 33336842
$ws.Range("K88").Value = 33336842
$ws.Range("M88").Value = -33336436

$ws.Range("H91").Value = 15194608
$ws.Range("I91").Value = 33336842
$ws.Range("K91").Value = 33336842
$ws.Range("M91").Value = -33335438

$ws.Range("H92").Value = 1365.4667
$ws.Range("I92").Value = 703.7
$ws.Range("K92").Value = 703.7
$ws.Range("M92").Value = 544.3

$ws.Range("H93").Value = 48997
$ws.Range("J93").Value = 48997
$ws.Range("L93").Value = 48997
$ws.Range("N93").Value = -53989

$ws.Range("H132").Value = 1647.5172
$ws.Range("I132").Value = 1670.6428
$ws.Range("K132").Value = 5011.928400000001
$ws.Range("M132").Value = -2481.928400000001

$ws.Range("H137").Value = 2681.2432
$ws.Range("I137").Value = 2413.4783
$ws.Range("K137").Value = 7240.4349
$ws.Range("M137").Value = -4690.4349

$ws.Range("H138").Value = 3837.9614
$ws.Range("I138").Value = 1217.1428
$ws.Range("J138").Value = 6895.5835
$ws.Range("K138").Value = 3651.4284
$ws.Range("L138").Value = 20686.7505
$ws.Range("M138").Value = 1488.5716
$ws.Range("N138").Value = -30966.7505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4676.1665
$ws.Range("I2").Value = 2329
$ws.Range("K2").Value = 2329
$ws.Range("M2").Value = -2216

$ws.Range("H45").Value = 9826.666999999999
$ws.Range("I45").Value = 2846.6
$ws.Range("K45").Value = 2846.6
$ws.Range("M45").Value = -2469.6

$ws.Range("H74").Value = 42150.73
$ws.Range("I74").Value = 81850.85000000001
$ws.Range("J74").Value = 4341.095
$ws.Range("K74").Value = 81850.85000000001
$ws.Range("L74").Value = 4341.095
$ws.Range("M74").Value = -80976.85000000001
$ws.Range("N74").Value = -6089.095

$ws.Range("H77").Value = 42150.73
$ws.Range("I77").Value = 81850.85000000001
$ws.Range("J77").Value = 4341.095
$ws.Range("K77").Value = 409254.25
$ws.Range("L77").Value = 21705.475
$ws.Range("M77").Value = -404886.25
$ws.Range("N77").Value = -30441.475

$ws.Range("H116").Value = 4676.1665
$ws.Range("I116").Value = 2329
$ws.Range("K116").Value = 2329
$ws.Range("M116").Value = -35

$ws.Range("H122").Value = 18380
$ws.Range("J122").Value = 6816.6665
$ws.Range("L122").Value = 20449.9995
$ws.Range("N122").Value = -25349.9995

$ws.Range("H141").Value = 74950
$ws.Range("J141").Value = 74950
$ws.Range("L141").Value = 74950
$ws.Range("N141").Value = -85310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4676.1665
$ws.Range("I3").Value = 2329
$ws.Range("K3").Value = 2329
$ws.Range("M3").Value = -2215

$ws.Range("H86").Value = 37077004
$ws.Range("I86").Value = 57676.61
$ws.Range("J86").Value = 111115656
$ws.Range("K86").Value = 57676.61
$ws.Range("L86").Value = 111115656
$ws.Range("M86").Value = -56553.61
$ws.Range("N86").Value = -111117902

$ws.Range("H89").Value = 37077004
$ws.Range("I89").Value = 57676.61
$ws.Range("J89").Value = 111115656
$ws.Range("K89").Value = 288383.05
$ws.Range("L89").Value = 555578280
$ws.Range("M89").Value = -282767.05
$ws.Range("N89").Value = -555589512

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 204.44444
$ws.Range("I22").Value = 173.33333
$ws.Range("K22").Value = 173.33333
$ws.Range("M22").Value = 176.66667

$ws.Range("H31").Value = 7587.974
$ws.Range("I31").Value = 3203.2727
$ws.Range("K31").Value = 3203.2727
$ws.Range("M31").Value = -2908.2727

$ws.Range("H34").Value = 7587.974
$ws.Range("I34").Value = 3203.2727
$ws.Range("K34").Value = 3203.2727
$ws.Range("M34").Value = -3001.2727

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H102").Value = 49241
$ws.Range("J102").Value = 49241
$ws.Range("L102").Value = 49241
$ws.Range("N102").Value = -54109

$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").ClearContents()

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H105").Value = 17865366
$ws.Range("I105").Value = 71428570
$ws.Range("J105").Value = 10964.333
$ws.Range("K105").Value = 71428570
$ws.Range("L105").Value = 10964.333
$ws.Range("M105").Value = -71426823
$ws.Range("N105").Value = -14458.333

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H132").Value = 6438.9707
$ws.Range("I132").Value = 2561.2856
$ws.Range("J132").Value = 9153.35
$ws.Range("K132").Value = 7683.8568
$ws.Range("L132").Value = 27460.05
$ws.Range("M132").Value = -5153.8568
$ws.Range("N132").Value = -32520.05

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3020.25
$ws.Range("I3").Value = 3020.25
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 9060.75
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -8948.75
$ws.Range("N3").ClearContents()

$ws.Range("H34").Value = 5901.231
$ws.Range("I34").Value = 261.5
$ws.Range("J34").Value = 6926.636
$ws.Range("K34").Value = 784.5
$ws.Range("L34").Value = 20779.908
$ws.Range("M34").Value = -700.5
$ws.Range("N34").Value = -20947.908

$ws.Range("H39").Value = 9607
$ws.Range("J39").Value = 11289.818
$ws.Range("L39").Value = 33869.454
$ws.Range("N39").Value = -34457.454

$ws.Range("H55").Value = 27089408
$ws.Range("I55").Value = 66666708
$ws.Range("J55").Value = 9099725
$ws.Range("K55").Value = 200000124
$ws.Range("L55").Value = 27299175
$ws.Range("M55").Value = -199999947
$ws.Range("N55").Value = -27299529

$ws.Range("H92").Value = 4049467
$ws.Range("J92").Value = 4525768.5
$ws.Range("L92").Value = 13577305.5
$ws.Range("N92").Value = -13579801.5

$ws.Range("H132").Value = 11251.762
$ws.Range("I132").Value = 5439.3
$ws.Range("J132").Value = 16535.818
$ws.Range("K132").Value = 48953.7
$ws.Range("L132").Value = 148822.362
$ws.Range("M132").Value = -46423.7
$ws.Range("N132").Value = -153882.362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 50000
$ws.Range("J34").Value = 50000
$ws.Range("L34").Value = 50000
$ws.Range("N34").Value = -50536

$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50630

$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52184

$ws.Range("H122").Value = 3295169.8
$ws.Range("J122").Value = 3750.7
$ws.Range("L122").Value = 11252.1
$ws.Range("N122").Value = -16152.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3936.4546
$ws.Range("I61").Value = 2062.348
$ws.Range("K61").Value = 2062.348
$ws.Range("M61").Value = -1860.348

$ws.Range("H82").Value = 1457.1666
$ws.Range("I82").Value = 899.75
$ws.Range("K82").Value = 899.75
$ws.Range("M82").Value = -538.75

$ws.Range("H85").Value = 1457.1666
$ws.Range("I85").Value = 899.75
$ws.Range("K85").Value = 899.75
$ws.Range("M85").Value = 348.25

$ws.Range("H113").Value = 3936.4546
$ws.Range("I113").Value = 2062.348
$ws.Range("K113").Value = 2062.348
$ws.Range("M113").Value = 107.652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 145184.39
$ws.Range("I122").Value = 211995
$ws.Range("K122").Value = 635985
$ws.Range("M122").Value = -633535

$ws.Range("H126").Value = 4587.636
$ws.Range("J126").Value = 5173
$ws.Range("L126").Value = 15519
$ws.Range("N126").Value = -20459
